# Update "想去人数" (want-to-go count) figures for several conan/anime
# event rows across the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 280
$wsExhibition.Range("F5").Value = 969
$wsExhibition.Range("F6").Value = 2316
$wsExhibition.Range("F7").Value = 199

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 280
$wsAll.Range("F7").Value = 969
$wsAll.Range("F8").Value = 2316
$wsAll.Range("F10").Value = 199
